# This workbook is a weekly price log. A new week's worth of data (2 rows,
# for the "Primera" quality split between "$/caja 36 atados" and
# "$/docena de atados" units) is inserted right before the existing row 759,
# pushing all the subsequent rows (old 759-873) down by two positions (to
# 761-875). The two freshly inserted rows are then populated with the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two whole rows at position 759 - everything at/after row 759 shifts
# down by two rows (759->761 ... 873->875), which is exactly the "shift"
# pattern seen across the rest of the sheet, and also naturally grows the
# used range from A1:R873 to A1:R875.
$ws.Range("759:760").Insert()

# Row 759: new "Primera" / "$/caja 36 atados" observation for the new date.
$ws.Range("A759").Value = 9
$ws.Range("B759").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C759").Value = "Metropolitana"
$ws.Range("D759").Value = 44984
$ws.Range("E759").Value = 13
$ws.Range("F759").Value = 100112040
$ws.Range("G759").Value = "Cilantro"
$ws.Range("H759").Value = "Sin especificar"
$ws.Range("I759").Value = "Primera"
$ws.Range("J759").Value = 34
$ws.Range("K759").Value = 10000
$ws.Range("L759").Value = 10000
$ws.Range("M759").Value = 10000
$ws.Range("N759").Value = "$/caja 36 atados"
$ws.Range("O759").Value = "Región Metropolitana"
$ws.Range("P759").Value = 278
$ws.Range("Q759").Value = 36
$ws.Range("R759").Value = "Hortaliza"

# Row 760: new "Primera" / "$/docena de atados" observation for the same date.
$ws.Range("A760").Value = 9
$ws.Range("B760").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C760").Value = "Metropolitana"
$ws.Range("D760").Value = 44984
$ws.Range("E760").Value = 13
$ws.Range("F760").Value = 100112040
$ws.Range("G760").Value = "Cilantro"
$ws.Range("H760").Value = "Sin especificar"
$ws.Range("I760").Value = "Primera"
$ws.Range("J760").Value = 70
$ws.Range("K760").Value = 18000
$ws.Range("L760").Value = 20000
$ws.Range("M760").Value = 19000
$ws.Range("N760").Value = "$/docena de atados"
$ws.Range("O760").Value = "Región Metropolitana"
$ws.Range("P760").Value = 6333
$ws.Range("Q760").Value = 3
$ws.Range("R760").Value = "Hortaliza"

Write-Output "ok"
